$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the daily team-pairing schedule (repartition des equipes par jour).
# "C A" was a bug; it should just be "C" (single team on duty), and the
# rotation pattern ("A B" / "C" / "B C") now also continues correctly
# through Samedi (row 16) and dimanche (row 17), which were previously
# left blank.

$ws.Range("C11").Value = "A B"
$ws.Range("D11").Value = "C"

$ws.Range("C12").Value = "B C"
$ws.Range("D12").Value = "A B"

$ws.Range("C13").Value = "C"
$ws.Range("D13").Value = "B C"

$ws.Range("C14").Value = "A B"
$ws.Range("D14").Value = "C"

$ws.Range("C15").Value = "B C"
$ws.Range("D15").Value = "A B"

$ws.Range("C16").Value = "C"
$ws.Range("D16").Value = "B C"

$ws.Range("C17").Value = "A B"
$ws.Range("D17").Value = "C"
